# 16.5.1.1b.xlsx — add a new "2022" data column (L) by extending the
# existing "2021" column (K): copy K2:K9 into L2:L9 so the new column
# inherits the same borders/fonts/number formats, then overwrite the
# copied values with the 2022 figures. Row 8's value additionally gets
# a thousands-separator number format, and the active selection moves
# to L2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (styles, borders, fonts) from column K into the new column L.
$src = $ws.Range("K2:K9")
$dst = $ws.Range("L2:L9")
$src.Copy($dst)

# Fill in the 2022 values.
$ws.Range("L3").Value = 2022
$ws.Range("L4").Value = 370
$ws.Range("L5").Value = 137
$ws.Range("L6").Value = 314
$ws.Range("L7").Value = 121
$ws.Range("L8").Value = 50
$ws.Range("L8").NumberFormat = "#,##0"
$ws.Range("L9").Value = 16

# Move the active selection to match the new cursor position.
$ws.Range("L2").Select()
